# "edit delete sort Program Module"
#
# Adds a new "AddProgram" worksheet (used for the module's Add/Edit/Delete/Sort
# tests) right after the existing "Program" sheet, seeds the "Program" sheet
# with a header row + one sample row, and seeds the new "AddProgram" sheet
# with its own header row + one sample row. The newly added sheet becomes the
# active tab.

$wb = $excel.ActiveWorkbook

# --- 1) Populate the existing "Program" sheet -----------------------------
# Header row, then the data row (Description/Status columns first, Name last)
# so shared-string indices line up the way the saved workbook expects.
$progSheet = $wb.Worksheets.Item("Program")
$progSheet.Range("A1").Value = "Name"
$progSheet.Range("B1").Value = "Description"
$progSheet.Range("C1").Value = "Status"
$progSheet.Range("B2").Value = "Java"
$progSheet.Range("C2").Value = "Active"
$progSheet.Range("A2").Value = "KarateSDET"
[void]$progSheet.Range("A2").Select()

# --- 2) Insert the new "AddProgram" sheet right after "Program" -----------
$afterSheet = $wb.Worksheets.Item("Program")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "AddProgram"

$newSheet.Range("A1").Value = "ProgramName"
$newSheet.Range("B1").Value = "ProgramDescription"
$newSheet.Range("C1").Value = "status"
$newSheet.Range("A2").Value = "FinalNinja"
$newSheet.Range("B2").Value = "Java"
$newSheet.Range("C2").Value = "Active"

# Make the new sheet the active tab/selection, as it was left after editing.
[void]$newSheet.Range("C2").Select()
